$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "use" column (G) from "Y" to "N" for the rows whose predictor
# type reclassification no longer marks them for use in the model.
$rowsToMarkUnused = @(12, 14, 17, 18, 19, 22, 24, 25, 26, 27, 29)
foreach ($r in $rowsToMarkUnused) {
    $ws.Cells.Item($r, 7).Value = "N"
}

# Update the "type" column (D) for row 26 ("range") from "boolean" to "num".
$ws.Cells.Item(26, 4).Value = "num"

# Reflect the final cell selection left active in the sheet.
$ws.Range("F14").Select()
